# Generate Report for Handoff
# Refreshes the "Latest Handoff Datetime" column for the rows whose handoff
# was (re-)generated in this run, on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$zhRows = @(7, 9, 10, 11, 12, 13, 14, 15, 16)
$deRows = @(7, 9, 10, 11, 12, 13, 14, 15, 16)

$zhTimestamp = "2016-03-08 14:22:10"
$deTimestamp = "2016-03-08 14:22:16"

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $zhRows) {
    $wsZh.Cells.Item($r, 4).Value = $zhTimestamp
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $deRows) {
    $wsDe.Cells.Item($r, 4).Value = $deTimestamp
}
